# Apply "adding support for port ranges" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Azure Non prod" -> "Azure Non prod VM range" and "Azure" -> "Azure Non prod VM range" in column A
$ws.Range("A2").Value = "Azure Non prod VM range"
$ws.Range("A3").Value = "Azure Non prod VM range"
$ws.Range("A4").Value = "Azure Non prod VM range"
$ws.Range("A5").Value = "Azure Non prod VM range"
$ws.Range("A6").Value = "Azure Non prod VM range"

# Rename "IOD- Non-Prod" -> "IOD- Non-Prod VM" and "Azure" -> "IOD- Non-Prod VM" in column C
$ws.Range("C2").Value = "IOD- Non-Prod VM"
$ws.Range("C3").Value = "IOD- Non-Prod VM"
$ws.Range("C4").Value = "IOD- Non-Prod VM"
$ws.Range("C5").Value = "IOD- Non-Prod VM"
$ws.Range("C6").Value = "IOD- Non-Prod VM"

# Update port value on row 2 to support a port range
$ws.Range("F2").Value = "22-25,  443"

# Add new rows 7-10 demonstrating udp/http port-range support
$ws.Range("A7").Value = "Azure Non prod"
$ws.Range("B7").Value = "10.108.0.1"
$ws.Range("C7").Value = "IOD- Non-Prod"
$ws.Range("D7").Value = "20.0.0.0/16"
$ws.Range("E7").Value = "udp"
$ws.Range("F7").Value = 80
$ws.Range("G7").Value = "http"
$ws.Range("H7").Value = "Overlap with rule at SerialNo:1"

$ws.Range("A8").Value = "Azure Non prod"
$ws.Range("B8").Value = "10.108.0.1"
$ws.Range("C8").Value = "IOD- Non-Prod"
$ws.Range("D8").Value = "20.0.0.0/16"
$ws.Range("E8").Value = "udp"
$ws.Range("F8").Value = 101
$ws.Range("G8").Value = "http"
$ws.Range("H8").Value = "Can be requested"

$ws.Range("A9").Value = "Azure Non prod"
$ws.Range("B9").Value = "10.108.0.1"
$ws.Range("C9").Value = "IOD- Non-Prod"
$ws.Range("D9").Value = "20.0.0.0/16"
$ws.Range("E9").Value = "udp"
$ws.Range("F9").Value = "80-105"
$ws.Range("G9").Value = "http"
$ws.Range("H9").Value = "Overlap with rule at SerialNo:1"

$ws.Range("A10").Value = "Azure Non prod"
$ws.Range("B10").Value = "10.108.0.1"
$ws.Range("C10").Value = "IOD- Non-Prod"
$ws.Range("D10").Value = "20.0.0.0/16"
$ws.Range("E10").Value = "udp"
$ws.Range("F10").Value = "101-105"
$ws.Range("G10").Value = "http"
$ws.Range("H10").Value = "Can be requested"
